$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 151 is a blank separator row (same pattern used between every day's
# block in this timesheet). Copy formatting from an existing separator
# row (140) so it gets the same shaded style as the others.
$ws.Range("A140:C140").Copy() | Out-Null
$ws.Range("A151:C151").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New timesheet entries for Feb 18 2020, appended starting row 152.
$entries = @(
    @("Feb 18 10:00 to 11:00", "Documented feature engineering code", $false),
    @("Feb 18 11:00 to 12:00", "Documented data cleaning and app code", $false),
    @("Feb 18 12:00 to 13:00", "worked on separate console message functionality, droped idea due to inconvinience`nin class decoration.", $true),
    @("Feb 18 13:00 to 13:30", "Working on data transformation, logic is little complicated", $false),
    @("Feb 18 13:30 to 14:00", "Lunch", $false),
    @("Feb 18 14:00 to 15:00", "Worked on data transformation, getting errors in process.", $false),
    @("Feb 18 15:00 to 16:00", "Logic failed for data transformation, using pycharm for debugging project", $false),
    @("Feb 18 16:00 to 17:00", "Modified some code for data transformation, data transformation working.", $false),
    @("Feb 18 17:00 to 18:00", "Modified code of progress bar, progress bars are working.", $false),
    @("Feb 18 18:00 to 19:00", "Created new features.", $false)
)

$startRow = 152
for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $entries[$i][0]
    $ws.Cells.Item($row, 2).Value = $entries[$i][1]
    $ws.Cells.Item($row, 3).Value = "Infimetrics"

    if ($entries[$i][2]) {
        # Long note wraps onto two lines, so Excel auto-expands the row.
        $ws.Cells.Item($row, 2).WrapText = $true
        $ws.Rows.Item($row).RowHeight = 30
    }
}

$lastRow = $startRow + $entries.Count - 1
$ws.Range("D$lastRow").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 143
